$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format so numeric-looking strings
# (prices, percentages) are preserved as text, matching the source data.
$ws.Range("D2:D50").NumberFormat = "@"
$ws.Range("E2:E50").NumberFormat = "@"

$ws.Range("D2").Value = "243.72"
$ws.Range("E2").Value = "-0.71%"
$ws.Range("D3").Value = "26.91"
$ws.Range("E3").Value = "4.02%"
$ws.Range("D4").Value = "5.155"
$ws.Range("E4").Value = "0.50%"
$ws.Range("D5").Value = "0.05617"
$ws.Range("E5").Value = "0.48%"
$ws.Range("D6").Value = "6.495"
$ws.Range("E6").Value = "0.15%"
$ws.Range("D7").Value = "0.8162"
$ws.Range("E7").Value = "-0.12%"
$ws.Range("D8").Value = "0.8309"
$ws.Range("E8").Value = "-1.86%"
$ws.Range("D9").Value = "0.1329"
$ws.Range("E9").Value = "-0.99%"
$ws.Range("D10").Value = "0.06931"
$ws.Range("E10").Value = "-0.43%"
$ws.Range("D11").Value = "0.02898"
$ws.Range("E11").Value = "1.60%"
$ws.Range("D12").Value = "0.09375"
$ws.Range("E12").Value = "-0.30%"
$ws.Range("D13").Value = "0.001527"
$ws.Range("E13").Value = "1.20%"
$ws.Range("D14").Value = "0.0005980"
$ws.Range("E14").Value = "-0.35%"
$ws.Range("D15").Value = "0.006175"
$ws.Range("E15").Value = "0.95%"
$ws.Range("D16").Value = "3.636"
$ws.Range("E16").Value = "2.87%"
$ws.Range("E17").Value = "-0.30%"
$ws.Range("D18").Value = "2.302"
$ws.Range("E18").Value = "8.69%"
$ws.Range("E19").Value = "-1.79%"
$ws.Range("E20").Value = "-4.61%"
$ws.Range("D21").Value = "0.1291"
$ws.Range("E21").Value = "-2.16%"
$ws.Range("D22").Value = "3.751"
$ws.Range("E22").Value = "0.33%"
$ws.Range("D23").Value = "0.04594"
$ws.Range("E23").Value = "-2.18%"
$ws.Range("E24").Value = "-2.43%"
$ws.Range("D25").Value = "0.001225"
$ws.Range("E25").Value = "-1.76%"
$ws.Range("D26").Value = "0.004487"
$ws.Range("E26").Value = "-2.58%"
$ws.Range("D27").Value = "0.00009800"
$ws.Range("E28").Value = "0.66%"
$ws.Range("D40").Value = "0.03638"
$ws.Range("E40").Value = "-0.49%"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006010"
$ws.Range("E41").Value = "-2.13%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1050"
$ws.Range("E42").Value = "-0.20%"
$ws.Range("D43").Value = "0.002618"
$ws.Range("E43").Value = "5.65%"
$ws.Range("D44").Value = "0.008952"
$ws.Range("E44").Value = "21.15%"
$ws.Range("D45").Value = "0.00005313"
$ws.Range("E45").Value = "-0.18%"
$ws.Range("E46").Value = "-0.02%"
$ws.Range("E47").Value = "-18.37%"
$ws.Range("D48").Value = "0.002609"
$ws.Range("E48").Value = "22.62%"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("E50").Value = "-0.02%"
